$d = $word.ActiveDocument

# Locate the word "joinctee", which sits inside its own <bp>...</bp>
# markup: " ou une <ms>petite <bp>joinctee</bp></ms>. Puys ..."
$rngJoin = $d.Content
$found = $rngJoin.Find.Execute("joinctee", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $joinStart = $rngJoin.Start
    $joinEnd = $rngJoin.End

    # "<bp>" immediately precedes "joinctee"
    $rngOpenBp = $d.Range($joinStart - 4, $joinStart)
    # "</bp>" immediately follows "joinctee"
    $rngCloseBp = $d.Range($joinEnd, $joinEnd + 5)

    if ($rngOpenBp.Text -eq "<bp>" -and $rngCloseBp.Text -eq "</bp>") {
        # Delete the closing tag first so the opening tag's offsets stay valid.
        $rngCloseBp.Delete()
        $rngOpenBp.Delete()
    }
}
